$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Men moet beschikken over Visual Studio Ultimate 2013 of Visual Studio Professional 2013 en GitHub.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Men moet beschikken over Windows 8.1 Professional, omdat we gebruik moeten maken van Microsoft Hyper-V. Verder moet men beschikken over Visual Studio Ultimate 2013 of Visual Studio Professional 2013 en GitHub.",
    2
)
